# =====================================================================
# Edit: add "2022-Q4" sheet with fund holdings data, insert as new
# second sheet (right after "总计"), and prepend a corresponding
# "2022-Q4" summary row to the "总计" sheet (shifting the other rows
# down by one).
# =====================================================================

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) Update "总计" (summary) sheet: insert 2022-Q4 as the new first
#    data row, shifting 2022-Q3 / 2021-Q4 / 2021-Q3 / 2020-Q4 down.
# ---------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

# Row 6 (2020-Q4) is brand new territory for this sheet (it used to
# stop at row 5), so it has no style yet. Clone the existing index
# column's style (bold/bordered "s=2") from A2 onto A6 before writing
# values, so every index cell A2:A6 ends up consistently styled.
$summary.Range("A2").Copy($summary.Range("A6"))

$summary.Cells.Item(2,1).Value = 0
$summary.Cells.Item(2,2).Value = "2022-Q4"
$summary.Cells.Item(2,3).Value = 15
$summary.Cells.Item(2,4).Value = 3.77

$summary.Cells.Item(3,1).Value = 1
$summary.Cells.Item(3,2).Value = "2022-Q3"
$summary.Cells.Item(3,3).Value = 19
$summary.Cells.Item(3,4).Value = 2.12

$summary.Cells.Item(4,1).Value = 2
$summary.Cells.Item(4,2).Value = "2021-Q4"
$summary.Cells.Item(4,3).Value = 1
$summary.Cells.Item(4,4).Value = 0.03

$summary.Cells.Item(5,1).Value = 3
$summary.Cells.Item(5,2).Value = "2021-Q3"
$summary.Cells.Item(5,3).Value = 2
$summary.Cells.Item(5,4).Value = 0.05

$summary.Cells.Item(6,1).Value = 4
$summary.Cells.Item(6,2).Value = "2020-Q4"
$summary.Cells.Item(6,3).Value = 2
$summary.Cells.Item(6,4).Value = 0.01

# ---------------------------------------------------------------
# 2) Insert the new "2022-Q4" worksheet right before "2022-Q3" (i.e.
#    as the second tab, right after "总计").
# ---------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($wb.Worksheets.Item("2022-Q3"))
$newSheet.Name = "2022-Q4"

# IMPORTANT: worksheet handles in this host track tab *position*, not
# a stable object identity, so re-resolve "2022-Q3" by name now that
# it has been pushed one slot to the right by the Add() above -- using
# the handle obtained before Add() would silently alias the brand new
# (still-empty) sheet instead.
$oldQ3 = $wb.Worksheets.Item("2022-Q3")

# Re-use the "2022-Q3" sheet's formatting as a template: its header
# row (B1:H1) and its index column (A) both use the bold/bordered
# "s=2" style, while the other data cells use the default style.
$oldQ3.Range("B1:H1").Copy($newSheet.Range("B1:H1"))
$oldQ3.Range("A2").Copy($newSheet.Range("A2:A16"))

# Columns B-G hold numeric-looking values that must stay TEXT (as in
# the source data), so force the number format to Text before writing.
$newSheet.Range("B2:G16").NumberFormat = "@"


$newSheet.Cells.Item(2,1).Value = 0
$newSheet.Cells.Item(2,2).Value = "001832"
$newSheet.Cells.Item(2,3).Value = "易方达瑞恒灵活配置混合"
$newSheet.Cells.Item(2,4).Value = "37.41"
$newSheet.Cells.Item(2,5).Value = "85.21"
$newSheet.Cells.Item(2,6).Value = "4.29"
$newSheet.Cells.Item(2,7).Value = "1.6049"
$newSheet.Cells.Item(2,8).Value = 9

$newSheet.Cells.Item(3,1).Value = 1
$newSheet.Cells.Item(3,2).Value = "160106"
$newSheet.Cells.Item(3,3).Value = "南方高增长混合（LOF）"
$newSheet.Cells.Item(3,4).Value = "16.27"
$newSheet.Cells.Item(3,5).Value = "88.50"
$newSheet.Cells.Item(3,6).Value = "6.93"
$newSheet.Cells.Item(3,7).Value = "1.1275"
$newSheet.Cells.Item(3,8).Value = 3

$newSheet.Cells.Item(4,1).Value = 2
$newSheet.Cells.Item(4,2).Value = "160105"
$newSheet.Cells.Item(4,3).Value = "南方积极配置混合（LOF）"
$newSheet.Cells.Item(4,4).Value = "5.52"
$newSheet.Cells.Item(4,5).Value = "89.92"
$newSheet.Cells.Item(4,6).Value = "6.78"
$newSheet.Cells.Item(4,7).Value = "0.3743"
$newSheet.Cells.Item(4,8).Value = 3

$newSheet.Cells.Item(5,1).Value = 3
$newSheet.Cells.Item(5,2).Value = "009234"
$newSheet.Cells.Item(5,3).Value = "鹏华优质企业混合"
$newSheet.Cells.Item(5,4).Value = "3.21"
$newSheet.Cells.Item(5,5).Value = "80.54"
$newSheet.Cells.Item(5,6).Value = "3.63"
$newSheet.Cells.Item(5,7).Value = "0.1165"
$newSheet.Cells.Item(5,8).Value = 7

$newSheet.Cells.Item(6,1).Value = 4
$newSheet.Cells.Item(6,2).Value = "002819"
$newSheet.Cells.Item(6,3).Value = "招商丰美灵活配置混合A"
$newSheet.Cells.Item(6,4).Value = "5.43"
$newSheet.Cells.Item(6,5).Value = "39.25"
$newSheet.Cells.Item(6,6).Value = "1.96"
$newSheet.Cells.Item(6,7).Value = "0.1064"
$newSheet.Cells.Item(6,8).Value = 10

$newSheet.Cells.Item(7,1).Value = 5
$newSheet.Cells.Item(7,2).Value = "014202"
$newSheet.Cells.Item(7,3).Value = "天弘中证1000指数增强C"
$newSheet.Cells.Item(7,4).Value = "6.60"
$newSheet.Cells.Item(7,5).Value = "94.11"
$newSheet.Cells.Item(7,6).Value = "1.59"
$newSheet.Cells.Item(7,7).Value = "0.1049"
$newSheet.Cells.Item(7,8).Value = 5

$newSheet.Cells.Item(8,1).Value = 6
$newSheet.Cells.Item(8,2).Value = "000554"
$newSheet.Cells.Item(8,3).Value = "南方中国梦灵活配置混合"
$newSheet.Cells.Item(8,4).Value = "1.50"
$newSheet.Cells.Item(8,5).Value = "92.96"
$newSheet.Cells.Item(8,6).Value = "6.82"
$newSheet.Cells.Item(8,7).Value = "0.1023"
$newSheet.Cells.Item(8,8).Value = 3

$newSheet.Cells.Item(9,1).Value = 7
$newSheet.Cells.Item(9,2).Value = "001753"
$newSheet.Cells.Item(9,3).Value = "红土创新新兴产业灵活配置混合"
$newSheet.Cells.Item(9,4).Value = "2.64"
$newSheet.Cells.Item(9,5).Value = "66.35"
$newSheet.Cells.Item(9,6).Value = "2.49"
$newSheet.Cells.Item(9,7).Value = "0.0657"
$newSheet.Cells.Item(9,8).Value = 9

$newSheet.Cells.Item(10,1).Value = 8
$newSheet.Cells.Item(10,2).Value = "014201"
$newSheet.Cells.Item(10,3).Value = "天弘中证1000指数增强A"
$newSheet.Cells.Item(10,4).Value = "3.86"
$newSheet.Cells.Item(10,5).Value = "94.11"
$newSheet.Cells.Item(10,6).Value = "1.59"
$newSheet.Cells.Item(10,7).Value = "0.0614"
$newSheet.Cells.Item(10,8).Value = 5

$newSheet.Cells.Item(11,1).Value = 9
$newSheet.Cells.Item(11,2).Value = "002389"
$newSheet.Cells.Item(11,3).Value = "招商安德灵活配置混合A"
$newSheet.Cells.Item(11,4).Value = "3.08"
$newSheet.Cells.Item(11,5).Value = "42.54"
$newSheet.Cells.Item(11,6).Value = "1.61"
$newSheet.Cells.Item(11,7).Value = "0.0496"
$newSheet.Cells.Item(11,8).Value = 10

$newSheet.Cells.Item(12,1).Value = 10
$newSheet.Cells.Item(12,2).Value = "002390"
$newSheet.Cells.Item(12,3).Value = "招商安德灵活配置混合C"
$newSheet.Cells.Item(12,4).Value = "2.06"
$newSheet.Cells.Item(12,5).Value = "42.54"
$newSheet.Cells.Item(12,6).Value = "1.61"
$newSheet.Cells.Item(12,7).Value = "0.0332"
$newSheet.Cells.Item(12,8).Value = 10

$newSheet.Cells.Item(13,1).Value = 11
$newSheet.Cells.Item(13,2).Value = "560006"
$newSheet.Cells.Item(13,3).Value = "益民核心增长混合"
$newSheet.Cells.Item(13,4).Value = "0.49"
$newSheet.Cells.Item(13,5).Value = "77.58"
$newSheet.Cells.Item(13,6).Value = "3.13"
$newSheet.Cells.Item(13,7).Value = "0.0153"
$newSheet.Cells.Item(13,8).Value = 3

$newSheet.Cells.Item(14,1).Value = 12
$newSheet.Cells.Item(14,2).Value = "015466"
$newSheet.Cells.Item(14,3).Value = "太平中证1000指数增强A"
$newSheet.Cells.Item(14,4).Value = "0.62"
$newSheet.Cells.Item(14,5).Value = "93.58"
$newSheet.Cells.Item(14,6).Value = "0.90"
$newSheet.Cells.Item(14,7).Value = "0.0056"
$newSheet.Cells.Item(14,8).Value = 9

$newSheet.Cells.Item(15,1).Value = 13
$newSheet.Cells.Item(15,2).Value = "002820"
$newSheet.Cells.Item(15,3).Value = "招商丰美灵活配置混合C"
$newSheet.Cells.Item(15,4).Value = "0.14"
$newSheet.Cells.Item(15,5).Value = "39.25"
$newSheet.Cells.Item(15,6).Value = "1.96"
$newSheet.Cells.Item(15,7).Value = "0.0027"
$newSheet.Cells.Item(15,8).Value = 10

$newSheet.Cells.Item(16,1).Value = 14
$newSheet.Cells.Item(16,2).Value = "015467"
$newSheet.Cells.Item(16,3).Value = "太平中证1000指数增强C"
$newSheet.Cells.Item(16,4).Value = "0.08"
$newSheet.Cells.Item(16,5).Value = "93.58"
$newSheet.Cells.Item(16,6).Value = "0.90"
$newSheet.Cells.Item(16,7).Value = "0.0007"
$newSheet.Cells.Item(16,8).Value = 9


Write-Host "2022-Q4 sheet added and 总计 sheet updated."
